$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, shifting existing data down
$ws.Rows.Item(1).Insert()

# Set the new header cell value
$ws.Range("A1").Value = "name file"

# Update the selection to match target state
$ws.Range("H7").Select()
